$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.235.39"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.893.84"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'243.32"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("E6").Value = "  +5.59%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "'41.40"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.347"
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("D10").Value = "'50.34"
$ws.Range("E10").Value = "  +8.08%  "
$ws.Range("D11").Value = "'0.0708"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "2.168.58"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("E14").Value = "  +5.01%  "
$ws.Range("D15").Value = "'0.692"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "1.881.70"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "'4.82"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "35.206.34"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'71.15"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "0.0₃0810"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "'240.85"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'12.41"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").Value = "'4.73"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  +32.79%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'169.76"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'8.38"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").Value = "'18.22"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "'4.11"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.939"
$ws.Range("E32").Value = "  +17.21%  "
$ws.Range("D33").Value = "'0.0559"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "'4.10"
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "'1.09"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("D41").Value = "'0.0637"
$ws.Range("E41").Value = "  +15.12%  "
$ws.Range("D42").Value = "'15.95"
$ws.Range("E42").Value = "  +7.79%  "
$ws.Range("D43").Value = "'88.99"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "1.337.71"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'48.52"
$ws.Range("E45").Value = "  +41.23%  "
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("D47").Value = "'2.41"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").Value = "'6.51"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "2.078.15"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").Value = "'11.40"
$ws.Range("E51").Value = "  -11.40%  "
